# Generate Report for Handback
#
# The localization pipeline re-processed the "2f65c9f5-f37b-45d5-973d-
# 36e67bd6949b.md" source file: a fresh xliff hand-off/hand-back round
# completed for both target languages, so the corresponding datetime
# stamps on the zh-cn and de-de report sheets move forward, and the
# Overview sheet's "Latest HO Xliff Generate Date" (max of the two
# per-language hand-off dates) is refreshed to match.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# zh-cn: row 2 is the 2f65c9f5-...md file.
#   H = Correspond Handoff Datetime, K = Correspond Handback DateTime
$wsZhCn.Range("H2").Value = "2016-08-13 02:58:49"
$wsZhCn.Range("K2").Value = "2016-08-13 02:59:17"

# de-de: row 2 is the 2f65c9f5-...md file.
#   H = Correspond Handoff Datetime, K = Correspond Handback DateTime
$wsDeDe.Range("H2").Value = "2016-08-13 02:58:56"
$wsDeDe.Range("K2").Value = "2016-08-13 02:59:26"

# Overview: row 2 is the 2f65c9f5-...md file.
#   G = Latest HO Xliff Generate Date (max of the per-language hand-off dates above)
$wsOverview.Range("G2").Value = "2016-08-13 02:58:56"
